# Add a new Job Posting row (JD_002 / Senior Engineer / sjdbsh / 1 / 4)
# underneath the existing JD_001 row, per commit message:
# "Add Job Posting with Job_Id=JD_002"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "JD_002"
$ws.Range("B3").Value = "Senior Engineer"
$ws.Range("C3").Value = "sjdbsh"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 4
